# Append a new, effectively-empty data row (row 2) below the header row
# on the active sheet ("Đơn sale chính"), expanding the used range from
# A1:T1 to A1:T2, mirroring an export that appended a blank record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold numeric totals/amounts end up as literal 0 in the
# appended blank row (I, K, L, M, N, O, P).
$zeroCols = @(9, 11, 12, 13, 14, 15, 16)
foreach ($col in $zeroCols) {
    $ws.Cells.Item(2, $col).Value = 0
}

# All remaining columns (A-T minus the numeric ones above) are blank
# text/empty cells in the appended row. Force Excel to keep a (blank)
# cell record at each of these addresses so the row is fully populated
# from A2 through T2, matching the exported row shape.
$blankCols = @(1, 2, 3, 4, 5, 6, 7, 8, 10, 17, 18, 19, 20)
foreach ($col in $blankCols) {
    $cell = $ws.Cells.Item(2, $col)
    $cell.NumberFormat = "General"
    $cell.Value = ""
}
